$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Un-merge the two title bands so we can shrink them from A:D to A:C
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").UnMerge()
$ws.Range("A10:D10").UnMerge()

# ---------------------------------------------------------------------------
# 2) Table 1 (rows 1-8) : "Conjunto de Teste (36 Amostras)"
#    - update title text
#    - update MAPE(%) values in column C
#    - remove the RMSE column (D) entirely
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Conjunto de Teste (36 Amostras)"

$ws.Range("C3").Value = 6392
$ws.Range("C4").Value = 5526
$ws.Range("C5").Value = 7128
$ws.Range("C6").Value = 5453
$ws.Range("C7").Value = 6232
$ws.Range("C8").Value = 5268

$ws.Range("D1:D8").Clear()
$ws.Range("D2").Clear()

# re-merge the title band to the new (narrower) extent
$ws.Range("A1:C1").Merge()

# ---------------------------------------------------------------------------
# 3) Table 2 (rows 10-17) : "Últimos 3 Meses (Última Amostra)"
#    - update title text
#    - update MAPE(%) values in column C
#    - remove the RMSE column (D) entirely, leaving behind the left-over
#      (unused / emptied) formatted cells that came from pasting the
#      cross-validation folds grid
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Últimos 3 Meses (Última Amostra)"

$ws.Range("C12").Value = 6734
$ws.Range("C13").Value = 2075
$ws.Range("C14").Value = 7015
$ws.Range("C15").Value = 1569
$ws.Range("C16").Value = 6872
$ws.Range("C17").Value = 1207

$ws.Range("D10:D17").Clear()

# re-merge the title band to the new (narrower) extent
$ws.Range("A10:C10").Merge()

# Left-over formatted (but empty) cells from the cross-validation paste,
# matching the ragged shape seen in the saved workbook.
$ws.Range("D12:D16").NumberFormat = "#,##0"

$ws.Range("E12:F12").NumberFormat = "#,##0"
$ws.Range("E12:F12").HorizontalAlignment = -4108
$ws.Range("G12:I12").NumberFormat = "#,##0"

$ws.Range("E13:F13").NumberFormat = "#,##0"
$ws.Range("E13:F13").HorizontalAlignment = -4108
$ws.Range("G13:H13").NumberFormat = "#,##0"

$ws.Range("E14:F14").NumberFormat = "#,##0"
$ws.Range("E14:F14").HorizontalAlignment = -4108
$ws.Range("G14:H14").NumberFormat = "#,##0"

$ws.Range("E15:F15").NumberFormat = "#,##0"
$ws.Range("E15:F15").HorizontalAlignment = -4108
$ws.Range("G15:H15").NumberFormat = "#,##0"

$ws.Range("E16:F16").NumberFormat = "#,##0"
$ws.Range("E16:F16").HorizontalAlignment = -4108
$ws.Range("G16:J16").NumberFormat = "#,##0"

$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("F17").HorizontalAlignment = -4108
$ws.Range("G17:H17").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# 4) Final selection, matching the state the workbook was saved in
# ---------------------------------------------------------------------------
$ws.Range("D17:H17").Select()
